$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Systems_EK1")
$ws.Activate()

# Update the hyperlink cell B2: shorten the displayed URL text and drop the
# "#/home" sub-address/location so it just points at the base hyperlink target.
$hlink = $ws.Hyperlinks.Item(1)
$hlink.TextToDisplay = "https://www.maut.toll-collect.de"
$hlink.SubAddress = ""

# Move the sheet's active selection from C5 to B2.
$ws.Range("B2").Select()
